$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The evaluation period "Davan Martinho do Nascimento" (row 20) had no
# recorded answers at all, so the whole row is removed from the weekly
# history table. Excel shifts every following row up by one and also
# keeps the SUM() formulas in the totals row in sync automatically.
$ws.Rows(20).Delete()

# The totals row (now row 29, used to be row 30) computes each column's
# percentage as SUM(col2:col-1-last-data-row)*100/<headcount>. With one
# fewer evaluation in the history, the headcount drops from 28 to 27.
# (Columns that never had any occurrence keep dividing by 28, matching
# the source workbook exactly.)
$colsToRebase = @("D","E","F","G","I","K","L","N","O","S","T","U","W","X","Y","AA","AB","AC","AD","AE","AF","AG","AH","AI")
foreach ($col in $colsToRebase) {
    $cell = $col + "29"
    $ws.Range($cell).Formula = "=SUM(" + $col + "2:" + $col + "28)*100/27"
}

# Restore the view so the selection lands back on the totals row, as in
# the authored workbook.
$ws.Range("F30").Select()
